# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# sheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# Sheet "展览" - column F updates (row -> new value)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 433
$ws1.Range("F5").Value = 1741
$ws1.Range("F6").Value = 88
$ws1.Range("F7").Value = 2187
$ws1.Range("F8").Value = 6
$ws1.Range("F11").Value = 4952
$ws1.Range("F21").Value = 3909
$ws1.Range("F23").Value = 677
$ws1.Range("F24").Value = 25
$ws1.Range("F31").Value = 582
$ws1.Range("F34").Value = 969
$ws1.Range("F35").Value = 2483

# Sheet "全部类型" - column F updates (row -> new value)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 433
$ws4.Range("F5").Value = 1741
$ws4.Range("F6").Value = 88
$ws4.Range("F7").Value = 2187
$ws4.Range("F8").Value = 6
$ws4.Range("F11").Value = 4952
$ws4.Range("F21").Value = 3909
$ws4.Range("F23").Value = 677
$ws4.Range("F24").Value = 25
$ws4.Range("F31").Value = 582
$ws4.Range("F35").Value = 969
$ws4.Range("F36").Value = 2483
